# spring 23 week 7 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E2").Value = 10.85

$ws.Range("D3").Value = 10.14
$ws.Range("F3").Value = 10.39

$ws.Range("C4").Value = 9.859999999999999
$ws.Range("E4").Value = 9.85
$ws.Range("F4").Value = 9.85

$ws.Range("B5").Value = 9.15
$ws.Range("D5").Value = 10.15
$ws.Range("F5").Value = 10.12
$ws.Range("H5").Value = 8.359999999999999
$ws.Range("I5").Value = 6.33

$ws.Range("C6").Value = 9.609999999999999
$ws.Range("D6").Value = 10.15
$ws.Range("E6").Value = 9.880000000000001
$ws.Range("G6").Value = 10.04

$ws.Range("F7").Value = 9.960000000000001
$ws.Range("H7").Value = 9.08

$ws.Range("E8").Value = 11.64
$ws.Range("G8").Value = 10.92

$ws.Range("E9").Value = 13.67
